$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("error message")

# Insert 7 new rows before existing row 32 (the "beas_mm_100x" material-management
# error rows used to live at the very end of the sheet; they get moved up here and
# a brand-new "beas_mm_1007" row is appended both here and at the very end).
$ws.Rows("32:38").Insert()

$newMsg = 'The "warehouse cost" and "consignment vendor" of "from warehouse" and "to warehouse" fields must be same.'

$rowsData = @(
    @(32, "beas_mm_1001", '<dw_1.item.itemcode.value>-<dw_1.item.itemname.value> cannot maintain in this screen'),
    @(33, "beas_mm_1002", '<cuser>,you have no right to open item master data window the screen will be close!'),
    @(34, "beas_mm_1003", "This item can NOT maintain in this screen"),
    @(35, "beas_mm_1004", "The Standard evaluation only"),
    @(36, "beas_mm_1005", "The unit name in Inventory Tab is mandatory information for inventory item"),
    @(37, "beas_mm_1006", "This item cost must great than 0, else it will be block in Transaction document")
)

foreach ($row in $rowsData) {
    $r = $row[0]
    $ws.Range("A$r").Value = "E   "
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("E$r").Value = $row[2]
}

# Row 38 introduces brand-new shared strings; the new message text (column E)
# was registered in the shared-string table before the new code (column C), so
# set E before C to reproduce the same shared-string ordering.
$ws.Range("A38").Value = "E   "
$ws.Range("E38").Value = $newMsg
$ws.Range("C38").Value = "beas_mm_1007"

# Append the same new "beas_mm_1007" row at the very end of the sheet (row 75).
$ws.Range("A75").Value = "E   "
$ws.Range("E75").Value = $newMsg
$ws.Range("C75").Value = "beas_mm_1007"

# Update the view's scroll/selection state to match the edited workbook.
$ws.Range("E35").Select()
$excel.ActiveWindow.ScrollRow = 20
